$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 23-25 (new Product/ProductType/PersonGender create API entries) ---
$ws.Range("B25").Value = "transaction.create.master.setProductType"
$ws.Range("B24").Value = "transaction.create.master.setProduct"
$ws.Range("B23").Value = "transaction.create.master.setPersonGender"

$ws.Range("C23").Value = "Menyimpan Data Baru Jenis Kelamin Orang"
$ws.Range("C24").Value = "Menyimpan Data Baru Produk"
$ws.Range("C25").Value = "Menyimpan Data Baru Jenis Produk"

# --- Insert 3 new rows before row 219 to host the update API entries ---
$ws.Range("216:218").Insert()

# The freshly inserted rows don't reliably inherit the B7/C4 row style from
# Insert() alone, so explicitly copy the formatting from the (now-shifted)
# row 219, which still carries the original B7/C4 styling.
$ws.Range("B219:C219").Copy()
$ws.Range("B216:C218").PasteSpecial(-4122)

$ws.Range("B218").Value = "transaction.update.master.setProductType"
$ws.Range("B217").Value = "transaction.update.master.setProduct"
$ws.Range("B216").Value = "transaction.update.master.setPersonGender"

$ws.Range("C216").Value = "Memutakhirkan Data Jenis Kelamin Orang"
$ws.Range("C217").Value = "Memutakhirkan Data Produk"
$ws.Range("C218").Value = "Memutakhirkan Data Jenis Produk"
